$d = $word.ActiveDocument

# The document currently ends with a "BodyText" styled paragraph
# ("BibTeX can be written like BibTeX."), immediately before the
# section properties. Insert a brand-new paragraph right after the
# very last paragraph in the document body.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
[void]$lastPara.Range.InsertParagraphAfter()

# The freshly-inserted (currently empty) paragraph is now the new
# last paragraph; grab its Range and replace its contents with the
# desired "BodyText" paragraph (style + run + text), expressed as raw
# WordprocessingML so the produced markup matches what Word itself
# would emit (including xml:space="preserve" on the text run).
$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newParaIndex)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$paragraphXml = "<w:p xmlns:w=`"$wNs`">" +
    "<w:pPr><w:pStyle w:val=`"BodyText`"/></w:pPr>" +
    "<w:r><w:t xml:space=`"preserve`">Various elipses are supported: ..., ..., ...</w:t></w:r>" +
    "</w:p>"

[void]$newPara.Range.InsertXML($paragraphXml)

Write-Output "Inserted ellipsis paragraph after paragraph $lastParaIndex."
